# insta_search_found.xlsx update:
#  - Row 6 (Benc4n) and Row 7 (tumblrindeed) lose their "Status" value (B6/B7 cleared).
#  - Row 8 (thesavagebean) keeps its "Status" value of OK.
#  - A new row 9 is appended for "tweeti_sai" with Status "OK".
#  - Selection/active cell moves to the newly added B9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the Status cells for Benc4n (row 6) and tumblrindeed (row 7).
$ws.Range("B6").ClearContents()
$ws.Range("B7").ClearContents()

# Append the new record in row 9.
$ws.Range("A9").Value = "tweeti_sai"
$ws.Range("B9").Value = "OK"

# Match the author's final selection: active cell B9, single-cell selection.
$ws.Range("B9").Select() | Out-Null
